$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6999.769
$ws.Range("I40").Value = 5714
$ws.Range("J40").Value = 8499.833000000001
$ws.Range("K40").Value = 5714
$ws.Range("L40").Value = 8499.833000000001
$ws.Range("M40").Value = -5539
$ws.Range("N40").Value = -8849.833000000001
$ws.Range("H64").Value = 8000
$ws.Range("J64").Value = 8000
$ws.Range("L64").Value = 8000
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 8000
$ws.Range("J67").Value = 8000
$ws.Range("L67").Value = 8000
$ws.Range("N67").Value = -9716
$ws.Range("H76").Value = 3460.8
$ws.Range("I76").Value = 3266.6667
$ws.Range("K76").Value = 3266.6667
$ws.Range("M76").Value = -2951.6667
$ws.Range("H79").Value = 3460.8
$ws.Range("I79").Value = 3266.6667
$ws.Range("K79").Value = 3266.6667
$ws.Range("M79").Value = -2174.6667
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H129").Value = 1475
$ws.Range("I129").Value = 1162.3334
$ws.Range("K129").Value = 3487.0002
$ws.Range("M129").Value = 1512.9998
$ws.Range("H137").Value = 2153.818
$ws.Range("I137").Value = 615.3333
$ws.Range("K137").Value = 1845.9999
$ws.Range("M137").Value = 704.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1081.6086
$ws.Range("I32").Value = 1030.4
$ws.Range("K32").Value = 1030.4
$ws.Range("M32").Value = -743.4000000000001
$ws.Range("H63").Value = 5776.8
$ws.Range("I63").Value = 2366
$ws.Range("K63").Value = 2366
$ws.Range("M63").Value = -1680
$ws.Range("H66").Value = 5776.8
$ws.Range("I66").Value = 2366
$ws.Range("K66").Value = 11830
$ws.Range("M66").Value = -8398
$ws.Range("H96").Value = 6688888.5
$ws.Range("J96").Value = 6688888.5
$ws.Range("L96").Value = 6688888.5
$ws.Range("N96").Value = -6694380.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 998
$ws.Range("I20").Value = 997.5
$ws.Range("K20").Value = 997.5
$ws.Range("M20").Value = -750.5
$ws.Range("H22").Value = 206
$ws.Range("I22").Value = 233.55556
$ws.Range("J22").Value = 82
$ws.Range("K22").Value = 233.55556
$ws.Range("L22").Value = 82
$ws.Range("M22").Value = -60.55556000000001
$ws.Range("N22").Value = -428
$ws.Range("H86").Value = 4246.933
$ws.Range("I86").Value = 1339.375
$ws.Range("K86").Value = 1339.375
$ws.Range("M86").Value = -216.375
$ws.Range("H89").Value = 4246.933
$ws.Range("I89").Value = 1339.375
$ws.Range("K89").Value = 6696.875
$ws.Range("M89").Value = -1080.875
$ws.Range("H105").Value = 1597.1538
$ws.Range("J105").Value = 1959.6666
$ws.Range("L105").Value = 1959.6666
$ws.Range("N105").Value = -5453.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 788.4286
$ws.Range("I22").Value = 536.8333
$ws.Range("K22").Value = 536.8333
$ws.Range("M22").Value = -186.8333
$ws.Range("H36").Value = 1366.6666
$ws.Range("I36").Value = 1366.6666
$ws.Range("K36").Value = 1366.6666
$ws.Range("M36").Value = -978.6666
$ws.Range("H40").Value = 1366.6666
$ws.Range("I40").Value = 1366.6666
$ws.Range("K40").Value = 1366.6666
$ws.Range("M40").Value = -1206.6666
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H68").Value = 66382.5
$ws.Range("J68").Value = 66382.5
$ws.Range("L68").Value = 66382.5
$ws.Range("N68").Value = -67880.5
$ws.Range("H71").Value = 66382.5
$ws.Range("J71").Value = 66382.5
$ws.Range("L71").Value = 199147.5
$ws.Range("N71").Value = -206635.5
$ws.Range("H96").Value = 7297.4
$ws.Range("J96").Value = 7297.4
$ws.Range("L96").Value = 7297.4
$ws.Range("N96").Value = -12789.4
$ws.Range("H107").Value = 708.93335
$ws.Range("I107").Value = 312
$ws.Range("K107").Value = 312
$ws.Range("M107").Value = 1608

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1564.8096
$ws.Range("I131").Value = 583.4167
$ws.Range("K131").Value = 1750.2501
$ws.Range("M131").Value = 3289.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12527
$ws.Range("J57").Value = 24999
$ws.Range("L57").Value = 24999
$ws.Range("N57").Value = -26639
$ws.Range("H80").Value = 4214.3335
$ws.Range("I80").Value = 4372.5
$ws.Range("J80").Value = 3898
$ws.Range("K80").Value = 4372.5
$ws.Range("L80").Value = 3898
$ws.Range("M80").Value = -3374.5
$ws.Range("N80").Value = -5894
$ws.Range("H83").Value = 4214.3335
$ws.Range("I83").Value = 4372.5
$ws.Range("J83").Value = 3898
$ws.Range("K83").Value = 21862.5
$ws.Range("L83").Value = 19490
$ws.Range("M83").Value = -16870.5
$ws.Range("N83").Value = -29474
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1547.091
$ws.Range("J22").Value = 1750
$ws.Range("L22").Value = 1750
$ws.Range("N22").Value = -2340
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H27").Value = 1547.091
$ws.Range("J27").Value = 1750
$ws.Range("L27").Value = 1750
$ws.Range("N27").Value = -1964
$ws.Range("H46").Value = 3161
$ws.Range("I46").Value = 2864.9167
$ws.Range("J46").Value = 3309.0417
$ws.Range("K46").Value = 2864.9167
$ws.Range("L46").Value = 3309.0417
$ws.Range("M46").Value = -2676.9167
$ws.Range("N46").Value = -3685.0417
$ws.Range("H68").Value = 4485.4287
$ws.Range("I68").Value = 3799.3333
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3799.3333
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -3050.3333
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 4485.4287
$ws.Range("I71").Value = 3799.3333
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 18996.6665
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -15252.6665
$ws.Range("N71").Value = -32488
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -21996
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 60000
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -69984
$ws.Range("H100").Value = 8183.0835
$ws.Range("I100").Value = 850
$ws.Range("J100").Value = 9649.700000000001
$ws.Range("K100").Value = 850
$ws.Range("L100").Value = 9649.700000000001
$ws.Range("M100").Value = -309
$ws.Range("N100").Value = -10731.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9100.299999999999
$ws.Range("I62").Value = 4751
$ws.Range("K62").Value = 4751
$ws.Range("M62").Value = -4127
$ws.Range("H65").Value = 9100.299999999999
$ws.Range("I65").Value = 4751
$ws.Range("K65").Value = 23755
$ws.Range("M65").Value = -20635
$ws.Range("H113").Value = 1078.6
$ws.Range("I113").Value = 1031.3334
$ws.Range("K113").Value = 3094.0002
$ws.Range("M113").Value = -924.0001999999999

